$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("alpha1F")

# Reuse the bold/bordered "id + label" style (cellXf index 1) for the three
# brand-new rows (17-19) by copying formats from existing styled rows, so
# Excel/IronCalc doesn't fabricate additional style entries.
$ws.Range("A13:M13").Copy()
$ws.Range("A17:M17").PasteSpecial(-4122)
$ws.Range("A14:M14").Copy()
$ws.Range("A18:M18").PasteSpecial(-4122)
$ws.Range("A15:M15").Copy()
$ws.Range("A19:M19").PasteSpecial(-4122)
$excel.CutCopyMode = 0

# --- Header row (row 2): same label text as before, re-set so every string
#     cell in the sheet is self-consistent after the table below is rebuilt. ---
$ws.Range("A2").Value = 0
$ws.Range("B2").Value = "HKL"
$ws.Range("C2").Value = "[1, 1, 0]"
$ws.Range("D2").Value = "[2, 0, 0]"
$ws.Range("E2").Value = "[2, 1, 1]"
$ws.Range("F2").Value = "[2, 2, 0]"
$ws.Range("G2").Value = "[3, 1, 0]"
$ws.Range("H2").Value = "[2, 2, 2]"
$ws.Range("I2").Value = "[3, 2, 1]"
$ws.Range("J2").Value = "[4, 0, 0]"
$ws.Range("K2").Value = "2Pairs"
$ws.Range("L2").Value = "4Pairs"
$ws.Range("M2").Value = "MaxUnique"

# --- Row 1 numeric column headers (unchanged) ---
$ws.Range("B1").Value = 0
$ws.Range("C1").Value = 1
$ws.Range("D1").Value = 2
$ws.Range("E1").Value = 3
$ws.Range("F1").Value = 4
$ws.Range("G1").Value = 5
$ws.Range("H1").Value = 6
$ws.Range("I1").Value = 7
$ws.Range("J1").Value = 8
$ws.Range("K1").Value = 9
$ws.Range("L1").Value = 10
$ws.Range("M1").Value = 11

# --- Data rows 3-9 (HKL 1-7): unchanged values, rewritten for completeness ---
# Row 3: ND Single
$ws.Range("A3").Value = 1
$ws.Range("B3").Value = "ND Single"
$ws.Range("C3").Value = 0.71
$ws.Range("D3").Value = 1.95
$ws.Range("E3").Value = 0.8
$ws.Range("F3").Value = 0.71
$ws.Range("G3").Value = 1.51
$ws.Range("H3").Value = 0.46
$ws.Range("I3").Value = 0.76
$ws.Range("J3").Value = 1.95
$ws.Range("K3").Value = 1.375
$ws.Range("L3").Value = 1.0425
$ws.Range("M3").Value = 1.031666666666667
# Row 4: RD Single
$ws.Range("A4").Value = 2
$ws.Range("B4").Value = "RD Single"
$ws.Range("C4").Value = 1.18
$ws.Range("D4").Value = 0.42
$ws.Range("E4").Value = 1.12
$ws.Range("F4").Value = 1.18
$ws.Range("G4").Value = 0.64
$ws.Range("H4").Value = 1.36
$ws.Range("I4").Value = 1.14
$ws.Range("J4").Value = 0.42
$ws.Range("K4").Value = 0.77
$ws.Range("L4").Value = 0.9749999999999999
$ws.Range("M4").Value = 0.9766666666666666
# Row 5: TD Single
$ws.Range("A5").Value = 3
$ws.Range("B5").Value = "TD Single"
$ws.Range("C5").Value = 1.38
$ws.Range("D5").Value = 0.52
$ws.Range("E5").Value = 0.99
$ws.Range("F5").Value = 1.38
$ws.Range("G5").Value = 0.78
$ws.Range("H5").Value = 0.97
$ws.Range("I5").Value = 1.09
$ws.Range("J5").Value = 0.52
$ws.Range("K5").Value = 0.755
$ws.Range("L5").Value = 1.0675
$ws.Range("M5").Value = 0.955
# Row 6: Morris
$ws.Range("A6").Value = 4
$ws.Range("B6").Value = "Morris"
$ws.Range("C6").Value = 0.99
$ws.Range("D6").Value = 0.35
$ws.Range("E6").Value = 1.25
$ws.Range("F6").Value = 0.99
$ws.Range("G6").Value = 0.58
$ws.Range("H6").Value = 1.74
$ws.Range("I6").Value = 1.17
$ws.Range("J6").Value = 0.35
$ws.Range("K6").Value = 0.8
$ws.Range("L6").Value = 0.895
$ws.Range("M6").Value = 1.013333333333333
# Row 7: Ring Perpendicular to ND
$ws.Range("A7").Value = 5
$ws.Range("B7").Value = "Ring Perpendicular to ND"
$ws.Range("C7").Value = 0.9335616438356165
$ws.Range("D7").Value = 1.31041095890411
$ws.Range("E7").Value = 0.9104109589041096
$ws.Range("F7").Value = 0.9335616438356165
$ws.Range("G7").Value = 1.179041095890411
$ws.Range("H7").Value = 0.7553424657534247
$ws.Range("I7").Value = 0.9134246575342466
$ws.Range("J7").Value = 1.31041095890411
$ws.Range("K7").Value = 1.11041095890411
$ws.Range("L7").Value = 1.021986301369863
$ws.Range("M7").Value = 1.000365296803653
# Row 8: Ring Perpendicular to RD
$ws.Range("A8").Value = 6
$ws.Range("B8").Value = "Ring Perpendicular to RD"
$ws.Range("C8").Value = 1.040526315789474
$ws.Range("D8").Value = 0.8184210526315789
$ws.Range("E8").Value = 1.04
$ws.Range("F8").Value = 1.040526315789474
$ws.Range("G8").Value = 0.891578947368421
$ws.Range("H8").Value = 1.12421052631579
$ws.Range("I8").Value = 1.041052631578947
$ws.Range("J8").Value = 0.8184210526315789
$ws.Range("K8").Value = 0.9292105263157895
$ws.Range("L8").Value = 0.9848684210526315
$ws.Range("M8").Value = 0.9926315789473685
# Row 9: Ring Perpendicular to TD
$ws.Range("A9").Value = 7
$ws.Range("B9").Value = "Ring Perpendicular to TD"
$ws.Range("C9").Value = 0.9510526315789474
$ws.Range("D9").Value = 0.9399999999999999
$ws.Range("E9").Value = 1.040526315789474
$ws.Range("F9").Value = 0.9510526315789474
$ws.Range("G9").Value = 0.9347368421052632
$ws.Range("H9").Value = 1.15421052631579
$ws.Range("I9").Value = 1.016842105263158
$ws.Range("J9").Value = 0.9399999999999999
$ws.Range("K9").Value = 0.9902631578947368
$ws.Range("L9").Value = 0.9706578947368421
$ws.Range("M9").Value = 1.006228070175439

# --- Data rows 10-19: new averaged-intensity table including the Spiral sampling
#     schemes (rows 11-13) ahead of the previously-existing NoRotation/Rotation/
#     HexGrid rows, which are now pushed down to rows 14-19. ---
# Row 10: Gaussian-Quadrature
$ws.Range("A10").Value = 8
$ws.Range("B10").Value = "Gaussian-Quadrature"
$ws.Range("C10").Value = 0.96339791054673
$ws.Range("D10").Value = 1.224158669892441
$ws.Range("E10").Value = 0.936824953861738
$ws.Range("F10").Value = 0.96339791054673
$ws.Range("G10").Value = 1.11310358006911
$ws.Range("H10").Value = 0.8511507835846838
$ws.Range("I10").Value = 0.9392363912951788
$ws.Range("J10").Value = 1.224158669892441
$ws.Range("K10").Value = 1.080491811877089
$ws.Range("L10").Value = 1.02194486121191
$ws.Range("M10").Value = 1.004645381541647
# Row 11: Spiral-90deg-10rot-5space
$ws.Range("A11").Value = 9
$ws.Range("B11").Value = "Spiral-90deg-10rot-5space"
$ws.Range("C11").Value = 0.9652252938229641
$ws.Range("D11").Value = 0.8715851201473213
$ws.Range("E11").Value = 1.058298695384647
$ws.Range("F11").Value = 0.9652252938229641
$ws.Range("G11").Value = 0.8950246625586119
$ws.Range("H11").Value = 1.204313551294188
$ws.Range("I11").Value = 1.034350441729079
$ws.Range("J11").Value = 0.8715851201473213
$ws.Range("K11").Value = 0.9649419077659841
$ws.Range("L11").Value = 0.9650836007944741
$ws.Range("M11").Value = 1.004799627489468
# Row 12: Spiral-90deg-15rot-5space
$ws.Range("A12").Value = 10
$ws.Range("B12").Value = "Spiral-90deg-15rot-5space"
$ws.Range("C12").Value = 0.9645701829220779
$ws.Range("D12").Value = 0.873008964621048
$ws.Range("E12").Value = 1.058049341446794
$ws.Range("F12").Value = 0.9645701829220779
$ws.Range("G12").Value = 0.8957705437731005
$ws.Range("H12").Value = 1.203813164357411
$ws.Range("I12").Value = 1.03402947822068
$ws.Range("J12").Value = 0.873008964621048
$ws.Range("K12").Value = 0.9655291530339207
$ws.Range("L12").Value = 0.9650496679779994
$ws.Range("M12").Value = 1.004873612556852
# Row 13: Spiral-90deg-10rot-3space
$ws.Range("A13").Value = 11
$ws.Range("B13").Value = "Spiral-90deg-10rot-3space"
$ws.Range("C13").Value = 0.9650113324059131
$ws.Range("D13").Value = 0.8720864908021738
$ws.Range("E13").Value = 1.058165888584316
$ws.Range("F13").Value = 0.9650113324059131
$ws.Range("G13").Value = 0.8952393132352363
$ws.Range("H13").Value = 1.20413249305916
$ws.Range("I13").Value = 1.034247895858543
$ws.Range("J13").Value = 0.8720864908021738
$ws.Range("K13").Value = 0.965126189693245
$ws.Range("L13").Value = 0.9650687610495791
$ws.Range("M13").Value = 1.004813902324224
# Row 14: NoRotation-tilt60deg
$ws.Range("A14").Value = 12
$ws.Range("B14").Value = "NoRotation-tilt60deg"
$ws.Range("C14").Value = 0.7804560000000001
$ws.Range("D14").Value = 1.568532
$ws.Range("E14").Value = 0.898532
$ws.Range("F14").Value = 0.7804560000000001
$ws.Range("G14").Value = 1.285431999999998
$ws.Range("H14").Value = 0.7516160000000014
$ws.Range("I14").Value = 0.8593719999999999
$ws.Range("J14").Value = 1.568532
$ws.Range("K14").Value = 1.233532
$ws.Range("L14").Value = 1.006994
$ws.Range("M14").Value = 1.02399
# Row 15: Rotation-NoTilt
$ws.Range("A15").Value = 13
$ws.Range("B15").Value = "Rotation-NoTilt"
$ws.Range("C15").Value = 0.71
$ws.Range("D15").Value = 1.95
$ws.Range("E15").Value = 0.798487499999998
$ws.Range("F15").Value = 0.71
$ws.Range("G15").Value = 1.509724999999999
$ws.Range("H15").Value = 0.46
$ws.Range("I15").Value = 0.76
$ws.Range("J15").Value = 1.95
$ws.Range("K15").Value = 1.374243749999999
$ws.Range("L15").Value = 1.042121874999999
$ws.Range("M15").Value = 1.03136875
# Row 16: Rotation-60detTilt
$ws.Range("A16").Value = 14
$ws.Range("B16").Value = "Rotation-60detTilt"
$ws.Range("C16").Value = 0.8280512774144025
$ws.Range("D16").Value = 1.551610225151999
$ws.Range("E16").Value = 0.8798463288320036
$ws.Range("F16").Value = 0.8280512774144025
$ws.Range("G16").Value = 1.293228532121598
$ws.Range("H16").Value = 0.6849694539776001
$ws.Range("I16").Value = 0.8599948802048017
$ws.Range("J16").Value = 1.551610225151999
$ws.Range("K16").Value = 1.215728276992001
$ws.Range("L16").Value = 1.021889777203202
$ws.Range("M16").Value = 1.016283449617067
# Row 17: HexGrid-90degTilt5degRes
$ws.Range("A17").Value = 15
$ws.Range("B17").Value = "HexGrid-90degTilt5degRes"
$ws.Range("C17").Value = 0.9968932899250197
$ws.Range("D17").Value = 0.9952835929948257
$ws.Range("E17").Value = 0.995247640063857
$ws.Range("F17").Value = 0.9968932899250197
$ws.Range("G17").Value = 0.9971977949902612
$ws.Range("H17").Value = 0.994016614433067
$ws.Range("I17").Value = 0.9957939515522974
$ws.Range("J17").Value = 0.9952835929948257
$ws.Range("K17").Value = 0.9952656165293414
$ws.Range("L17").Value = 0.9960794532271805
$ws.Range("M17").Value = 0.9957388139932214
# Row 18: HexGrid-90degTilt22p5degRes
$ws.Range("A18").Value = 16
$ws.Range("B18").Value = "HexGrid-90degTilt22p5degRes"
$ws.Range("C18").Value = 1.004582345526417
$ws.Range("D18").Value = 0.9656412496961175
$ws.Range("E18").Value = 1.001911299166746
$ws.Range("F18").Value = 1.004582345526417
$ws.Range("G18").Value = 0.9795921116313026
$ws.Range("H18").Value = 1.012626276966431
$ws.Range("I18").Value = 1.003090995685177
$ws.Range("J18").Value = 0.9656412496961175
$ws.Range("K18").Value = 0.9837762744314317
$ws.Range("L18").Value = 0.9941793099789245
$ws.Range("M18").Value = 0.9945740464453651
# Row 19: HexGrid-60degTilt5degRes
$ws.Range("A19").Value = 17
$ws.Range("B19").Value = "HexGrid-60degTilt5degRes"
$ws.Range("C19").Value = 1.017938947771571
$ws.Range("D19").Value = 0.898195602646797
$ws.Range("E19").Value = 1.018630548160817
$ws.Range("F19").Value = 1.017938947771571
$ws.Range("G19").Value = 0.9389302154134186
$ws.Range("H19").Value = 1.062578713043334
$ws.Range("I19").Value = 1.022262537378944
$ws.Range("J19").Value = 0.898195602646797
$ws.Range("K19").Value = 0.9584130754038069
$ws.Range("L19").Value = 0.9881760115876892
$ws.Range("M19").Value = 0.9930894274024804
